# Scheduled-runner update: refresh currentAveragePrice / Leve profit
# columns (H:N) across the per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM,
# LTW, WVR) with the latest market-board figures.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 618.567  # H15: 657.7217000000001 -> 618.567
$ws.Cells.Item(15, 9).Value = 618.567  # I15: 657.7217000000001 -> 618.567
$ws.Cells.Item(15, 11).Value = 1855.701  # K15: 1973.1651 -> 1855.701
$ws.Cells.Item(15, 13).Value = -1686.701  # M15: -1804.1651 -> -1686.701

$ws.Cells.Item(98, 8).Value = 1001.6957  # H98: 1028.6086 -> 1001.6957
$ws.Cells.Item(98, 9).Value = 843.06665  # I98: 884.3333 -> 843.06665
$ws.Cells.Item(98, 11).Value = 843.06665  # K98: 884.3333 -> 843.06665
$ws.Cells.Item(98, 13).Value = 654.93335  # M98: 613.6667 -> 654.93335

$ws.Cells.Item(113, 8).Value = 97675.734  # H113: 121260.836 -> 97675.734
$ws.Cells.Item(113, 9).Value = 112318.16  # I113: 160903.44 -> 112318.16
$ws.Cells.Item(113, 10).Value = 2500  # J113: 2333 -> 2500
$ws.Cells.Item(113, 11).Value = 112318.16  # K113: 160903.44 -> 112318.16
$ws.Cells.Item(113, 12).Value = 2500  # L113: 2333 -> 2500
$ws.Cells.Item(113, 13).Value = -109064.16  # M113: -157649.44 -> -109064.16
$ws.Cells.Item(113, 14).Value = -9008  # N113: -8841 -> -9008

$ws.Cells.Item(122, 8).Value = 1001.6957  # H122: 1028.6086 -> 1001.6957
$ws.Cells.Item(122, 9).Value = 843.06665  # I122: 884.3333 -> 843.06665
$ws.Cells.Item(122, 11).Value = 2529.19995  # K122: 2652.9999 -> 2529.19995
$ws.Cells.Item(122, 13).Value = -79.19995000000017  # M122: -202.9998999999998 -> -79.19995000000017

$ws.Cells.Item(125, 8).Value = 1916.6666  # H125: 1651.1111 -> 1916.6666
$ws.Cells.Item(125, 9).Value = 1916.6666  # I125: 2104 -> 1916.6666
$ws.Cells.Item(125, 10).Value = 0  # J125: 745.3333 -> 0
$ws.Cells.Item(125, 11).Value = 17249.9994  # K125: 18936 -> 17249.9994
$ws.Cells.Item(125, 12).Value = 0  # L125: 6707.9997 -> 0
$ws.Cells.Item(125, 13).Value = -14789.9994  # M125: -16476 -> -14789.9994
$ws.Cells.Item(125, 14).ClearContents()  # N125: cleared (was -11627.9997)

$ws.Cells.Item(132, 8).Value = 3574583.5  # H132: 4101940.8 -> 3574583.5
$ws.Cells.Item(132, 9).Value = 3201.8596  # I132: 3611.52 -> 3201.8596
$ws.Cells.Item(132, 10).Value = 19233718  # J132: 22730710 -> 19233718
$ws.Cells.Item(132, 11).Value = 9605.578799999999  # K132: 10834.56 -> 9605.578799999999
$ws.Cells.Item(132, 12).Value = 57701154  # L132: 68192130 -> 57701154
$ws.Cells.Item(132, 13).Value = -7075.578799999999  # M132: -8304.559999999999 -> -7075.578799999999
$ws.Cells.Item(132, 14).Value = -57706214  # N132: -68197190 -> -57706214

$ws.Cells.Item(137, 8).Value = 5406293  # H137: 7143793 -> 5406293
$ws.Cells.Item(137, 9).Value = 701.7826  # I137: 692.5 -> 701.7826
$ws.Cells.Item(137, 10).Value = 14286908  # J137: 14286893 -> 14286908
$ws.Cells.Item(137, 11).Value = 2105.3478  # K137: 2077.5 -> 2105.3478
$ws.Cells.Item(137, 12).Value = 42860724  # L137: 42860679 -> 42860724
$ws.Cells.Item(137, 13).Value = 444.6522  # M137: 472.5 -> 444.6522
$ws.Cells.Item(137, 14).Value = -42865824  # N137: -42865779 -> -42865824

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10217.976  # H32: 9501.955 -> 10217.976
$ws.Cells.Item(32, 9).Value = 9889.059999999999  # I32: 9022.311 -> 9889.059999999999
$ws.Cells.Item(32, 10).Value = 11595.3125  # J32: 11720.3125 -> 11595.3125
$ws.Cells.Item(32, 11).Value = 9889.059999999999  # K32: 9022.311 -> 9889.059999999999
$ws.Cells.Item(32, 12).Value = 11595.3125  # L32: 11720.3125 -> 11595.3125
$ws.Cells.Item(32, 13).Value = -9602.059999999999  # M32: -8735.311 -> -9602.059999999999
$ws.Cells.Item(32, 14).Value = -12169.3125  # N32: -12294.3125 -> -12169.3125

$ws.Cells.Item(61, 8).Value = 8773368  # H61: 14707878 -> 8773368
$ws.Cells.Item(61, 9).Value = 11112512  # I61: 16130843 -> 11112512
$ws.Cells.Item(61, 10).Value = 1579.0834  # J61: 3904.6667 -> 1579.0834
$ws.Cells.Item(61, 11).Value = 11112512  # K61: 16130843 -> 11112512
$ws.Cells.Item(61, 12).Value = 1579.0834  # L61: 3904.6667 -> 1579.0834
$ws.Cells.Item(61, 13).Value = -11112300  # M61: -16130631 -> -11112300
$ws.Cells.Item(61, 14).Value = -2003.0834  # N61: -4328.6667 -> -2003.0834

$ws.Cells.Item(74, 8).Value = 7814672  # H74: 7814660.5 -> 7814672
$ws.Cells.Item(74, 9).Value = 9435400  # I74: 9616843 -> 9435400
$ws.Cells.Item(74, 10).Value = 5712.909  # J74: 5203.5 -> 5712.909
$ws.Cells.Item(74, 11).Value = 9435400  # K74: 9616843 -> 9435400
$ws.Cells.Item(74, 12).Value = 5712.909  # L74: 5203.5 -> 5712.909
$ws.Cells.Item(74, 13).Value = -9434526  # M74: -9615969 -> -9434526
$ws.Cells.Item(74, 14).Value = -7460.909  # N74: -6951.5 -> -7460.909

$ws.Cells.Item(77, 8).Value = 7814672  # H77: 7814660.5 -> 7814672
$ws.Cells.Item(77, 9).Value = 9435400  # I77: 9616843 -> 9435400
$ws.Cells.Item(77, 10).Value = 5712.909  # J77: 5203.5 -> 5712.909
$ws.Cells.Item(77, 11).Value = 47177000  # K77: 48084215 -> 47177000
$ws.Cells.Item(77, 12).Value = 28564.545  # L77: 26017.5 -> 28564.545
$ws.Cells.Item(77, 13).Value = -47172632  # M77: -48079847 -> -47172632
$ws.Cells.Item(77, 14).Value = -37300.545  # N77: -34753.5 -> -37300.545

$ws.Cells.Item(97, 8).Value = 4661.8276  # H97: 7422.8887 -> 4661.8276
$ws.Cells.Item(97, 9).Value = 4911.8184  # I97: 8882.5 -> 4911.8184
$ws.Cells.Item(97, 10).Value = 3876.1428  # J97: 4503.6665 -> 3876.1428
$ws.Cells.Item(97, 11).Value = 4911.8184  # K97: 8882.5 -> 4911.8184
$ws.Cells.Item(97, 12).Value = 3876.1428  # L97: 4503.6665 -> 3876.1428
$ws.Cells.Item(97, 13).Value = -4415.8184  # M97: -8386.5 -> -4415.8184
$ws.Cells.Item(97, 14).Value = -4868.1428  # N97: -5495.6665 -> -4868.1428

$ws.Cells.Item(136, 8).Value = 8773368  # H136: 14707878 -> 8773368
$ws.Cells.Item(136, 9).Value = 11112512  # I136: 16130843 -> 11112512
$ws.Cells.Item(136, 10).Value = 1579.0834  # J136: 3904.6667 -> 1579.0834
$ws.Cells.Item(136, 11).Value = 33337536  # K136: 48392529 -> 33337536
$ws.Cells.Item(136, 12).Value = 4737.2502  # L136: 11714.0001 -> 4737.2502
$ws.Cells.Item(136, 13).Value = -33334986  # M136: -48389979 -> -33334986
$ws.Cells.Item(136, 14).Value = -9837.2502  # N136: -16814.0001 -> -9837.2502

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2266.4465  # H134: 2602.3022 -> 2266.4465
$ws.Cells.Item(134, 9).Value = 1359.3256  # I134: 1571.4166 -> 1359.3256
$ws.Cells.Item(134, 10).Value = 5266.923  # J134: 7904 -> 5266.923
$ws.Cells.Item(134, 11).Value = 4077.976799999999  # K134: 4714.2498 -> 4077.976799999999
$ws.Cells.Item(134, 12).Value = 15800.769  # L134: 23712 -> 15800.769
$ws.Cells.Item(134, 13).Value = -1542.976799999999  # M134: -2179.2498 -> -1542.976799999999
$ws.Cells.Item(134, 14).Value = -20870.769  # N134: -28782 -> -20870.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1206.2858  # H58: 1132.9778 -> 1206.2858
$ws.Cells.Item(58, 9).Value = 571.0714  # I58: 523.54285 -> 571.0714
$ws.Cells.Item(58, 10).Value = 2476.7144  # J58: 3266 -> 2476.7144
$ws.Cells.Item(58, 11).Value = 571.0714  # K58: 523.54285 -> 571.0714
$ws.Cells.Item(58, 12).Value = 2476.7144  # L58: 3266 -> 2476.7144
$ws.Cells.Item(58, 13).Value = -368.0714  # M58: -320.54285 -> -368.0714
$ws.Cells.Item(58, 14).Value = -2882.7144  # N58: -3672 -> -2882.7144

$ws.Cells.Item(99, 9).Value = 1200  # I99: 1150 -> 1200
$ws.Cells.Item(99, 10).Value = 2457  # J99: 2507 -> 2457
$ws.Cells.Item(99, 11).Value = 1200  # K99: 1150 -> 1200
$ws.Cells.Item(99, 12).Value = 2457  # L99: 2507 -> 2457
$ws.Cells.Item(99, 13).Value = 298  # M99: 348 -> 298
$ws.Cells.Item(99, 14).Value = -5453  # N99: -5503 -> -5453

$ws.Cells.Item(126, 9).Value = 1200  # I126: 1150 -> 1200
$ws.Cells.Item(126, 10).Value = 2457  # J126: 2507 -> 2457
$ws.Cells.Item(126, 11).Value = 3600  # K126: 3450 -> 3600
$ws.Cells.Item(126, 12).Value = 7371  # L126: 7521 -> 7371
$ws.Cells.Item(126, 13).Value = -1130  # M126: -980 -> -1130
$ws.Cells.Item(126, 14).Value = -12311  # N126: -12461 -> -12311

$ws.Cells.Item(132, 8).Value = 6494889.5  # H132: 7464234 -> 6494889.5
$ws.Cells.Item(132, 9).Value = 8065729  # I132: 9435336 -> 8065729
$ws.Cells.Item(132, 10).Value = 2086.2666  # J132: 2202.4285 -> 2086.2666
$ws.Cells.Item(132, 11).Value = 24197187  # K132: 28306008 -> 24197187
$ws.Cells.Item(132, 12).Value = 6258.7998  # L132: 6607.2855 -> 6258.7998
$ws.Cells.Item(132, 13).Value = -24194657  # M132: -28303478 -> -24194657
$ws.Cells.Item(132, 14).Value = -11318.7998  # N132: -11667.2855 -> -11318.7998

$ws.Cells.Item(134, 8).Value = 284710.6  # H134: 398369.84 -> 284710.6
$ws.Cells.Item(134, 9).Value = 1035.8308  # I134: 1294.5625 -> 1035.8308
$ws.Cells.Item(134, 10).Value = 1255176.9  # J134: 1986671 -> 1255176.9
$ws.Cells.Item(134, 11).Value = 3107.4924  # K134: 3883.6875 -> 3107.4924
$ws.Cells.Item(134, 12).Value = 3765530.7  # L134: 5960013 -> 3765530.7
$ws.Cells.Item(134, 13).Value = -572.4924000000001  # M134: -1348.6875 -> -572.4924000000001
$ws.Cells.Item(134, 14).Value = -3770600.7  # N134: -5965083 -> -3770600.7

$ws.Cells.Item(136, 8).Value = 1206.2858  # H136: 1132.9778 -> 1206.2858
$ws.Cells.Item(136, 9).Value = 571.0714  # I136: 523.54285 -> 571.0714
$ws.Cells.Item(136, 10).Value = 2476.7144  # J136: 3266 -> 2476.7144
$ws.Cells.Item(136, 11).Value = 1713.2142  # K136: 1570.62855 -> 1713.2142
$ws.Cells.Item(136, 12).Value = 7430.1432  # L136: 9798 -> 7430.1432
$ws.Cells.Item(136, 13).Value = 836.7857999999999  # M136: 979.3714499999999 -> 836.7857999999999
$ws.Cells.Item(136, 14).Value = -12530.1432  # N136: -14898 -> -12530.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(93, 8).Value = 8750  # H93: 6188.8887 -> 8750
$ws.Cells.Item(93, 10).Value = 8750  # J93: 6188.8887 -> 8750
$ws.Cells.Item(93, 12).Value = 26250  # L93: 18566.6661 -> 26250
$ws.Cells.Item(93, 14).Value = -29994  # N93: -22310.6661 -> -29994

$ws.Cells.Item(113, 8).Value = 718.575  # H113: 730.85895 -> 718.575
$ws.Cells.Item(113, 9).Value = 437.2  # I113: 438.06668 -> 437.2
$ws.Cells.Item(113, 10).Value = 1080.3429  # J113: 1130.1212 -> 1080.3429
$ws.Cells.Item(113, 11).Value = 1311.6  # K113: 1314.20004 -> 1311.6
$ws.Cells.Item(113, 12).Value = 3241.0287  # L113: 3390.3636 -> 3241.0287
$ws.Cells.Item(113, 13).Value = 858.4000000000001  # M113: 855.7999599999998 -> 858.4000000000001
$ws.Cells.Item(113, 14).Value = -7581.028700000001  # N113: -7730.363600000001 -> -7581.028700000001

$ws.Cells.Item(131, 8).Value = 806.48  # H131: 802.27 -> 806.48
$ws.Cells.Item(131, 9).Value = 451.15384  # I131: 453.92856 -> 451.15384
$ws.Cells.Item(131, 10).Value = 859.5747  # J131: 858.97675 -> 859.5747
$ws.Cells.Item(131, 11).Value = 1353.46152  # K131: 1361.78568 -> 1353.46152
$ws.Cells.Item(131, 12).Value = 2578.7241  # L131: 2576.93025 -> 2578.7241
$ws.Cells.Item(131, 13).Value = 3686.53848  # M131: 3678.21432 -> 3686.53848
$ws.Cells.Item(131, 14).Value = -12658.7241  # N131: -12656.93025 -> -12658.7241

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2257.805  # H132: 2232.035 -> 2257.805
$ws.Cells.Item(132, 9).Value = 1515.1077  # I132: 1547.7344 -> 1515.1077
$ws.Cells.Item(132, 10).Value = 5097.5293  # J132: 4222.727 -> 5097.5293
$ws.Cells.Item(132, 11).Value = 4545.3231  # K132: 4643.2032 -> 4545.3231
$ws.Cells.Item(132, 12).Value = 15292.5879  # L132: 12668.181 -> 15292.5879
$ws.Cells.Item(132, 13).Value = -2015.3231  # M132: -2113.2032 -> -2015.3231
$ws.Cells.Item(132, 14).Value = -20352.5879  # N132: -17728.181 -> -20352.5879

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5589.6787  # H7: 6920.6816 -> 5589.6787
$ws.Cells.Item(7, 9).Value = 7467.778  # I7: 12990.8 -> 7467.778
$ws.Cells.Item(7, 10).Value = 4700.0527  # J7: 5135.353 -> 4700.0527
$ws.Cells.Item(7, 11).Value = 7467.778  # K7: 12990.8 -> 7467.778
$ws.Cells.Item(7, 12).Value = 4700.0527  # L7: 5135.353 -> 4700.0527
$ws.Cells.Item(7, 13).Value = -7355.778  # M7: -12878.8 -> -7355.778
$ws.Cells.Item(7, 14).Value = -4924.0527  # N7: -5359.353 -> -4924.0527

$ws.Cells.Item(93, 8).Value = 1038.7307  # H93: 1350 -> 1038.7307
$ws.Cells.Item(93, 9).Value = 1033.5625  # I93: 0 -> 1033.5625
$ws.Cells.Item(93, 10).Value = 1047  # J93: 1350 -> 1047
$ws.Cells.Item(93, 11).Value = 1033.5625  # K93: 0 -> 1033.5625
$ws.Cells.Item(93, 12).Value = 1047  # L93: 1350 -> 1047
$ws.Cells.Item(93, 13).Value = 214.4375  # M93: None -> 214.4375
$ws.Cells.Item(93, 14).Value = -3543  # N93: -3846 -> -3543

$ws.Cells.Item(122, 8).Value = 4374.2246  # H122: 4251.827 -> 4374.2246
$ws.Cells.Item(122, 9).Value = 4088.5  # I122: 3967.353 -> 4088.5
$ws.Cells.Item(122, 10).Value = 4912.0586  # J122: 4789.1665 -> 4912.0586
$ws.Cells.Item(122, 11).Value = 12265.5  # K122: 11902.059 -> 12265.5
$ws.Cells.Item(122, 12).Value = 14736.1758  # L122: 14367.4995 -> 14736.1758
$ws.Cells.Item(122, 13).Value = -9815.5  # M122: -9452.059000000001 -> -9815.5
$ws.Cells.Item(122, 14).Value = -19636.1758  # N122: -19267.4995 -> -19636.1758

$ws.Cells.Item(126, 8).Value = 5589.6787  # H126: 6920.6816 -> 5589.6787
$ws.Cells.Item(126, 9).Value = 7467.778  # I126: 12990.8 -> 7467.778
$ws.Cells.Item(126, 10).Value = 4700.0527  # J126: 5135.353 -> 4700.0527
$ws.Cells.Item(126, 11).Value = 22403.334  # K126: 38972.39999999999 -> 22403.334
$ws.Cells.Item(126, 12).Value = 14100.1581  # L126: 15406.059 -> 14100.1581
$ws.Cells.Item(126, 13).Value = -19933.334  # M126: -36502.39999999999 -> -19933.334
$ws.Cells.Item(126, 14).Value = -19040.1581  # N126: -20346.059 -> -19040.1581

$ws.Cells.Item(132, 8).Value = 3227.28  # H132: 5887766 -> 3227.28
$ws.Cells.Item(132, 9).Value = 2406.5432  # I132: 3249.7576 -> 2406.5432
$ws.Cells.Item(132, 10).Value = 6726.2104  # J132: 26328716 -> 6726.2104
$ws.Cells.Item(132, 11).Value = 7219.6296  # K132: 9749.272799999999 -> 7219.6296
$ws.Cells.Item(132, 12).Value = 20178.6312  # L132: 78986148 -> 20178.6312
$ws.Cells.Item(132, 13).Value = -4689.6296  # M132: -7219.272799999999 -> -4689.6296
$ws.Cells.Item(132, 14).Value = -25238.6312  # N132: -78991208 -> -25238.6312

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1312.0834  # H132: 1397.418 -> 1312.0834
$ws.Cells.Item(132, 9).Value = 1057.836  # I132: 1090.2373 -> 1057.836
$ws.Cells.Item(132, 10).Value = 2722  # J132: 3662.875 -> 2722
$ws.Cells.Item(132, 11).Value = 3173.508  # K132: 3270.7119 -> 3173.508
$ws.Cells.Item(132, 12).Value = 8166  # L132: 10988.625 -> 8166
$ws.Cells.Item(132, 13).Value = -643.5079999999998  # M132: -740.7119000000002 -> -643.5079999999998
$ws.Cells.Item(132, 14).Value = -13226  # N132: -16048.625 -> -13226

$ws.Cells.Item(136, 8).Value = 856.2320999999999  # H136: 950.25 -> 856.2320999999999
$ws.Cells.Item(136, 9).Value = 655  # I136: 726.63416 -> 655
$ws.Cells.Item(136, 10).Value = 1781.9  # J136: 2260 -> 1781.9
$ws.Cells.Item(136, 11).Value = 1965  # K136: 2179.90248 -> 1965
$ws.Cells.Item(136, 12).Value = 5345.700000000001  # L136: 6780 -> 5345.700000000001
$ws.Cells.Item(136, 13).Value = 585  # M136: 370.0975200000003 -> 585
$ws.Cells.Item(136, 14).Value = -10445.7  # N136: -11880 -> -10445.7
